# Apply odds updates to the FlashScore sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Odd_Over2_FT (Q3) and Odd_Under2_FT (R3): 1.93 -> 1.87
$ws.Range("Q3").Value = 1.87
$ws.Range("R3").Value = 1.87

# Row 4: Odd_Over15_FT (O4): 1.29 -> 1.25
$ws.Range("O4").Value = 1.25
# Row 4: Odd_Under15_FT (P4): 3.5 -> 3.75
$ws.Range("P4").Value = 3.75
# Row 4: Odd_Over25_FT (S4): 1.9 -> 1.82
$ws.Range("S4").Value = 1.82
# Row 4: Odd_Under25_FT (T4): 1.95 -> 1.92
$ws.Range("T4").Value = 1.92
